$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update C5 value
$ws.Range("C5").Value = "1:30pm-2:15pm 9/2/2021"

# 2. Insert a new row at row 11 (shifts existing rows 11+ down by one)
$ws.Rows("11").Insert()

# 3. Populate the new row 11 with data
$ws.Range("A11").Value = "Melbourne"
$ws.Range("B11").Value = "Brunetti: Terminal 4, Melbourne Airport"
$ws.Range("C11").Value = "4:45am - 1:15pm, 9/2/2021"
$ws.Range("D11").Value = "Case attended venue"
